# Updated cryptos list on Sat Feb 17 22:17:21 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a "Price" (column D) cell while preventing
# Excel's COM automation from auto-coercing plain-decimal-looking strings
# (e.g. "352.93") into real numbers. Forcing the cell to Text format first
# makes the assignment stick as a literal string, matching how these sheets
# store prices as text.
function Set-TextValue($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.716.12"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.783.21"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "352.93"
$ws.Range("E5").Value = "  -1.54%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "108.97"
$ws.Range("E6").Value = "  -0.21%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.550"
$ws.Range("E7").Value = "  -2.89%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.598"
$ws.Range("E9").Value = "  +0.74%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "39.87"
$ws.Range("E10").Value = "  -0.51%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +2.63%  "

# Row 12 - now Dogecoin (was Chainlink)
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D12") "0.0837"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13 - now Chainlink (was Dogecoin)
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "20.11"
$ws.Range("E13").Value = "  +3.04%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.65"
$ws.Range("E14").Value = "  +0.65%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.221.84"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.794.03"
$ws.Range("E16").Value = "  +0.82%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.927"
$ws.Range("E17").Value = "  -2.30%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "51.710.67"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "7.77"
$ws.Range("E19").Value = "  +4.85%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("D20") "3.13"
$ws.Range("E20").Value = "  -0.16%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("E21").Value = "  +0.94%  "

# Row 22 - ShibaInu (price has a subscript "3" character, U+2083)
$sub3 = [char]0x2083
$ws.Range("D22").Value = "0.0{0}0965" -f $sub3
$ws.Range("E22").Value = "  -1.72%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "69.87"
$ws.Range("E23").Value = "  -0.54%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "266.64"
$ws.Range("E24").Value = "  -2.74%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.60%  "

# Row 26 - now Dai (was EthereumClassic)
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27 - now EthereumClassic (was Dai)
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "26.12"
$ws.Range("E27").Value = "  -2.11%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +11.91%  "

# Row 29 - Cosmos
Set-TextValue $ws.Range("D29") "10.22"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30 - now Toncoin (was InjectiveProtocol)
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.23"
$ws.Range("E30").Value = "  -2.48%  "

# Row 31 - now InjectiveProtocol (was Filecoin)
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D31") "36.73"
$ws.Range("E31").Value = "  +6.43%  "

# Row 32 - now Filecoin (was OKB)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "6.25"
$ws.Range("E32").Value = "  +9.33%  "

# Row 33 - now OKB (was VeChain)
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D33") "51.83"
$ws.Range("E33").Value = "  +0.54%  "

# Row 34 - now VeChain (was RenderToken)
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D34") "0.0453"
$ws.Range("E34").Value = "  -2.46%  "

# Row 35 - now RenderToken (was Toncoin)
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D35") "5.54"
$ws.Range("E35").Value = "  +4.84%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.0832"
$ws.Range("E36").Value = "  -1.47%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.06%  "

# Row 38 - Celestia
$ws.Range("E38").Value = "  +2.87%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -2.68%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.86%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -0.17%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -0.80%  "

# Row 43 - Monero
Set-TextValue $ws.Range("D43") "120.73"
$ws.Range("E43").Value = "  -0.95%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "22.06"
$ws.Range("E44").Value = "  +0.19%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -2.83%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.127.91"
$ws.Range("E46").Value = "  +2.61%  "

# Row 47 - NEARProtocol
Set-TextValue $ws.Range("D47") "3.30"
$ws.Range("E47").Value = "  +1.56%  "

# Row 48 - ApeXProtocol
Set-TextValue $ws.Range("D48") "2.33"
$ws.Range("E48").Value = "  +7.13%  "

# Row 49 - now SEI (was THORChain)
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
Set-TextValue $ws.Range("D49") "0.908"
$ws.Range("E49").Value = "  -2.43%  "

# Row 50 - now THORChain (was SEI)
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "5.43"
$ws.Range("E50").Value = "  -5.31%  "

# Row 51 - TrustWalletToken
Set-TextValue $ws.Range("D51") "1.33"
$ws.Range("E51").Value = "  +8.76%  "
